# Update "想去人数" (want-to-go count) values in column F across sheets.
# Mirrors the upstream data refresh (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1665
$ws1.Range("F6").Value = 620
$ws1.Range("F7").Value = 1107
$ws1.Range("F8").Value = 1550
$ws1.Range("F12").Value = 1461
$ws1.Range("F13").Value = 3085
$ws1.Range("F14").Value = 624
$ws1.Range("F15").Value = 1768
$ws1.Range("F16").Value = 1804
$ws1.Range("F17").Value = 857
$ws1.Range("F18").Value = 279
$ws1.Range("F20").Value = 1473
$ws1.Range("F23").Value = 10
$ws1.Range("F24").Value = 1214
$ws1.Range("F25").Value = 404
$ws1.Range("F26").Value = 458
$ws1.Range("F27").Value = 115
$ws1.Range("F28").Value = 4793
$ws1.Range("F29").Value = 36
$ws1.Range("F30").Value = 751
$ws1.Range("F31").Value = 573
$ws1.Range("F32").Value = 1656
$ws1.Range("F33").Value = 72
$ws1.Range("F34").Value = 139

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 65
$ws2.Range("F6").Value = 56
$ws2.Range("F7").Value = 74

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 65
$ws4.Range("F9").Value = 56
$ws4.Range("F10").Value = 74
$ws4.Range("F12").Value = 1665
$ws4.Range("F14").Value = 620
$ws4.Range("F15").Value = 1107
$ws4.Range("F16").Value = 1550
$ws4.Range("F21").Value = 1461
$ws4.Range("F22").Value = 3085
$ws4.Range("F23").Value = 624
$ws4.Range("F24").Value = 1768
$ws4.Range("F25").Value = 1804
$ws4.Range("F26").Value = 857
$ws4.Range("F27").Value = 279
$ws4.Range("F29").Value = 1473
$ws4.Range("F33").Value = 10
$ws4.Range("F35").Value = 1214
$ws4.Range("F36").Value = 404
$ws4.Range("F37").Value = 458
$ws4.Range("F38").Value = 115
$ws4.Range("F39").Value = 4793
$ws4.Range("F40").Value = 36
$ws4.Range("F41").Value = 751
$ws4.Range("F42").Value = 573
$ws4.Range("F43").Value = 1656
$ws4.Range("F46").Value = 72
$ws4.Range("F47").Value = 139
